# Update for new ebewe dates and compliance periods
# Shifts the compliance schedule forward: each building-ID bucket now shows
# its next cycle of dates (Initial Compliance Due Date, Initial Comparative
# Period, and Next Compliance Due Date), and the old asterisked / rich-text
# "Sept 7, 2023*" footnote values are replaced with plain dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "0 or 1"
$ws.Range("B2").Value = "Dec 1, 2026"
$ws.Range("C2").Value = "Dec 1, 2011 - Dec 1, 2026"
$ws.Range("D2").Value = "Dec 1, 2031"

# Row 3: "2 or 3"
$ws.Range("B3").Value = "Dec 1, 2027"
$ws.Range("C3").Value = "Dec 1, 2022 - Dec 1, 2027"
$ws.Range("D3").Value = "Dec 1, 2032"

# Row 4: "4 or 5"
$ws.Range("B4").Value = "Dec 1, 2028"
$ws.Range("C4").Value = "Dec 1, 2023 - Dec 1, 2028"
$ws.Range("D4").Value = "Dec 1, 2033"

# Row 5: "6 or 7"
$ws.Range("B5").Value = "Dec 1, 2029"
$ws.Range("C5").Value = "Dec 1, 2024 - Dec 1, 2029"
$ws.Range("D5").Value = "Dec 1, 2034"

# Row 6: "8 or 9" (unchanged values, but now formatted as Text like the rest
# of the column so the whole column shares one consistent number format)
$ws.Range("B6").Value = "Dec 1, 2025"
$ws.Range("C6").Value = "Dec 1, 2020 - Dec 1, 2025"
$ws.Range("D6").Value = "Dec 1, 2030"

# Column B (Initial Compliance Due Date) is now stored as Text for every
# data row, matching the format already used elsewhere in the sheet.
$ws.Range("B2:B6").NumberFormat = "@"

# Rows 2 & 3 no longer contain wrapped multi-line rich text (the bold
# "Sept 7, 2023*" footnote run is gone), so let the rows shrink back to the
# sheet's default height instead of staying pinned at the old taller size.
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()

# Leave the active selection where the author left off editing.
$ws.Range("C6").Select() | Out-Null
